# The "treatment of metal scrap, mixed, for recycling, unsorted, sorting"
# activity's contribution values (rows 7 and 8, columns A:P of the contribution
# table) had the wrong sign; flip the sign of the numeric value embedded in
# each of those cells' text while leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
$rows = @(7, 8)

foreach ($col in $columns) {
    foreach ($row in $rows) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $v = $cell.Value2

        if ($v -ne $null -and $v.ToString().Contains("treatment of metal scrap")) {
            $text = $v.ToString()

            if ($text -match '^(.*,\s*)(-?[0-9.eE+\-]+)(\s*\]\]\s*)$') {
                $prefix = $matches[1]
                $numStr = $matches[2]
                $suffix = $matches[3]

                if ($numStr.StartsWith('-')) {
                    $negStr = $numStr.Substring(1)
                }
                else {
                    $negStr = '-' + $numStr
                }

                $newVal = $prefix + $negStr + $suffix
                $cell.Value = $newVal
            }
        }
    }
}
